# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.365.44"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.627.22"
$ws.Range("D4").Value = "'0.9997"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'0.9996"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'302.55"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("D7").Value = "'0.3757"
$ws.Range("E7").Value = "  +0.83%  "
$ws.Range("D8").Value = "'0.3624"
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "'51.38"
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("E10").Value = "  +0.64%  "
$ws.Range("D11").Value = "'1.217"
$ws.Range("E11").Value = "  -2.46%  "
$ws.Range("D12").Value = "'1.0000"
$ws.Range("D13").Value = "'22.20"
$ws.Range("E13").Value = "  -2.63%  "
$ws.Range("D14").Value = "'6.459"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").Value = "'0.00001237"
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("D16").Value = "'7.289"
$ws.Range("E16").Value = "  +0.13%  "
$ws.Range("D17").Value = "1.602.45"
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("D18").Value = "'94.46"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").Value = "'0.06933"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("D20").Value = "'17.54"
$ws.Range("E20").Value = "  -3.08%  "
$ws.Range("D21").Value = "'6.546"
$ws.Range("E21").Value = "  +0.66%  "
$ws.Range("D22").Value = "'1.0000"
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").Value = "23.363.72"
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("D25").Value = "'2.473"
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("D26").Value = "'3.055"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").Value = "'21.11"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").Value = "'149.89"
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").Value = "'5.269"
$ws.Range("E29").Value = "  -1.13%  "
$ws.Range("D30").Value = "'132.52"
$ws.Range("E30").Value = "  -2.41%  "
$ws.Range("D31").Value = "1.787.05"
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").Value = "'6.586"
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("D33").Value = "'2.151"
$ws.Range("E33").Value = "  -5.72%  "
$ws.Range("D34").Value = "'1.056"
$ws.Range("E34").Value = "  +11.24%  "
$ws.Range("D35").Value = "'11.15"
$ws.Range("E35").Value = "  +7.79%  "
$ws.Range("D36").Value = "'0.02752"
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("D37").Value = "'0.08757"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "'0.2486"
$ws.Range("E38").Value = "  -1.39%  "
$ws.Range("D39").Value = "'0.07116"
$ws.Range("E39").Value = "  -1.19%  "
$ws.Range("D40").Value = "'5.968"
$ws.Range("E40").Value = "  -2.31%  "
$ws.Range("D41").Value = "'0.6970"
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("D42").Value = "'1.327"
$ws.Range("E42").Value = "  -3.14%  "
$ws.Range("D43").Value = "'15.75"
$ws.Range("E43").Value = "  -1.63%  "
$ws.Range("D44").Value = "'12.00"
$ws.Range("E44").Value = "  -3.52%  "
$ws.Range("D45").Value = "'0.6435"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").Value = "'0.9991"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "'2.270"
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("D48").Value = "'3.956"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "'0.07970"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "'127.15"
$ws.Range("E50").Value = "  -0.86%  "
$ws.Range("D51").Value = "'1.187"
$ws.Range("E51").Value = "  -0.82%  "
